$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.563.74"
$ws.Range("E2").Value = "  +5.80%  "
$ws.Range("D3").Value = "1.709.22"
$ws.Range("E3").Value = "  +4.33%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "222.69"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").Value = "0.536"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "29.88"
$ws.Range("E8").Value = "  +3.85%  "
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("D10").Value = "0.0651"
$ws.Range("E10").Value = "  +7.08%  "
$ws.Range("D11").Value = "0.0910"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "1.955.02"
$ws.Range("E12").Value = "  +4.38%  "
$ws.Range("D13").Value = "1.700.20"
$ws.Range("E13").Value = "  +4.20%  "
$ws.Range("D14").Value = "0.614"
$ws.Range("E14").Value = "  +4.41%  "
$ws.Range("D15").Value = "10.15"
$ws.Range("E15").Value = "  +7.84%  "
$ws.Range("E16").Value = "  +8.23%  "
$ws.Range("D17").Value = "31.544.53"
$ws.Range("E17").Value = "  +5.69%  "
$ws.Range("D18").Value = "67.37"
$ws.Range("E18").Value = "  +4.71%  "
$ws.Range("D19").Value = "251.09"
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "10.14"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "159.40"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "16.07"
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "3.91"
$ws.Range("E30").Value = "  +15.76%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  +4.24%  "
$ws.Range("D33").Value = "3.41"
$ws.Range("E33").Value = "  +6.94%  "
$ws.Range("D34").Value = "1.533.47"
$ws.Range("E34").Value = "  +8.06%  "
$ws.Range("D35").Value = "1.75"
$ws.Range("E35").Value = "  +4.13%  "
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").Value = "0.614"
$ws.Range("E37").Value = "  +9.37%  "
$ws.Range("D38").Value = "82.64"
$ws.Range("E38").Value = "  +8.50%  "
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("D41").Value = "2.31"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "2.05"
$ws.Range("E42").Value = "  +5.60%  "
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("D44").Value = "0.0504"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "1.04"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "52.39"
$ws.Range("E47").Value = "  +5.61%  "
$ws.Range("D48").Value = "5.61"
$ws.Range("E48").Value = "  +5.04%  "
$ws.Range("D49").Value = "1.847.86"
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("E50").Value = "  +9.49%  "
$ws.Range("D51").Value = "93.60"
$ws.Range("E51").Value = "  +0.68%  "
